# Apply the edits described by the diff:
# 1. Update the date paragraph from 2025-06-13 to 2025-06-16
# 2. Delete the "Timeline: URL link" paragraph within the
#    "Guide: Public Art - Sprengel Museum" section
# 3. Update the folium map memory address for that same section

$d = $word.ActiveDocument

# 1. Update date
$d.Content.Find.Execute("2025-06-13", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-06-16", 2)

# 2. Remove the "Timeline: URL link" paragraph (first occurrence only —
#    the Guide: Public Art - Sprengel Museum section)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Timeline: URL link") {
        $p.Range.Delete()
        break
    }
}

# 3. Update the folium map hex address for that section
$d.Content.Find.Execute("0x23bede4d2a0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0x1d7981b0460", 2)
